$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 204601.4
$ws.Range("I107").Value = 253252
$ws.Range("K107").Value = 253252
$ws.Range("M107").Value = -251332
$ws.Range("H112").Value = 1811.4117
$ws.Range("J112").Value = 1811.4117
$ws.Range("L112").Value = 5434.2351
$ws.Range("N112").Value = -7650.2351
$ws.Range("H132").Value = 3876.7827
$ws.Range("I132").Value = 4114.0527
$ws.Range("K132").Value = 12342.1581
$ws.Range("M132").Value = -9812.158100000001
$ws.Range("H134").Value = 70824.234
$ws.Range("J134").Value = 70824.234
$ws.Range("L134").Value = 70824.234
$ws.Range("N134").Value = -80964.234
$ws.Range("H138").Value = 6284.803
$ws.Range("I138").Value = 3322.2173
$ws.Range("K138").Value = 9966.651899999999
$ws.Range("M138").Value = -4826.651899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4563.8057
$ws.Range("I61").Value = 4485.793
$ws.Range("J61").Value = 4887
$ws.Range("K61").Value = 4485.793
$ws.Range("L61").Value = 4887
$ws.Range("M61").Value = -4273.793
$ws.Range("N61").Value = -5311
$ws.Range("H74").Value = 2485.4
$ws.Range("I74").Value = 2726.5386
$ws.Range("K74").Value = 2726.5386
$ws.Range("M74").Value = -1852.5386
$ws.Range("H77").Value = 2485.4
$ws.Range("I77").Value = 2726.5386
$ws.Range("K77").Value = 13632.693
$ws.Range("M77").Value = -9264.692999999999
$ws.Range("H132").Value = 3470.7144
$ws.Range("I132").Value = 2584.6
$ws.Range("J132").Value = 5686
$ws.Range("K132").Value = 7753.799999999999
$ws.Range("L132").Value = 17058
$ws.Range("M132").Value = -5223.799999999999
$ws.Range("N132").Value = -22118
$ws.Range("H133").Value = 65000
$ws.Range("J133").Value = 65000
$ws.Range("L133").Value = 65000
$ws.Range("N133").Value = -70060
$ws.Range("H134").Value = 49800
$ws.Range("J134").Value = 49800
$ws.Range("L134").Value = 49800
$ws.Range("N134").Value = -59940
$ws.Range("H136").Value = 4563.8057
$ws.Range("I136").Value = 4485.793
$ws.Range("J136").Value = 4887
$ws.Range("K136").Value = 13457.379
$ws.Range("L136").Value = 14661
$ws.Range("M136").Value = -10907.379
$ws.Range("N136").Value = -19761

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 163.22223
$ws.Range("I80").Value = 81.5
$ws.Range("J80").Value = 228.6
$ws.Range("K80").Value = 81.5
$ws.Range("L80").Value = 228.6
$ws.Range("M80").Value = 916.5
$ws.Range("N80").Value = -2224.6
$ws.Range("H83").Value = 163.22223
$ws.Range("I83").Value = 81.5
$ws.Range("J83").Value = 228.6
$ws.Range("K83").Value = 407.5
$ws.Range("L83").Value = 1143
$ws.Range("M83").Value = 4584.5
$ws.Range("N83").Value = -11127
$ws.Range("H105").Value = 2160.56
$ws.Range("I105").Value = 2181.2273
$ws.Range("K105").Value = 2181.2273
$ws.Range("M105").Value = -434.2273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3990.7812
$ws.Range("I16").Value = 3541
$ws.Range("J16").Value = 5340.125
$ws.Range("K16").Value = 3541
$ws.Range("L16").Value = 5340.125
$ws.Range("M16").Value = -3254
$ws.Range("N16").Value = -5914.125
$ws.Range("H31").Value = 43824.082
$ws.Range("I31").Value = 1433.3334
$ws.Range("J31").Value = 57954.332
$ws.Range("K31").Value = 1433.3334
$ws.Range("L31").Value = 57954.332
$ws.Range("M31").Value = -1138.3334
$ws.Range("N31").Value = -58544.332
$ws.Range("H34").Value = 43824.082
$ws.Range("I34").Value = 1433.3334
$ws.Range("J34").Value = 57954.332
$ws.Range("K34").Value = 1433.3334
$ws.Range("L34").Value = 57954.332
$ws.Range("M34").Value = -1231.3334
$ws.Range("N34").Value = -58358.332
$ws.Range("H58").Value = 2608.35
$ws.Range("I58").Value = 2780.125
$ws.Range("K58").Value = 2780.125
$ws.Range("M58").Value = -2577.125
$ws.Range("H113").Value = 3990.7812
$ws.Range("I113").Value = 3541
$ws.Range("J113").Value = 5340.125
$ws.Range("K113").Value = 3541
$ws.Range("L113").Value = 5340.125
$ws.Range("M113").Value = -1371
$ws.Range("N113").Value = -9680.125
$ws.Range("H136").Value = 2608.35
$ws.Range("I136").Value = 2780.125
$ws.Range("K136").Value = 8340.375
$ws.Range("M136").Value = -5790.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 5568.4614
$ws.Range("I56").Value = 5568.4614
$ws.Range("K56").Value = 5568.4614
$ws.Range("M56").Value = -5038.4614
$ws.Range("H93").Value = 696.5
$ws.Range("I93").Value = 696.5
$ws.Range("K93").Value = 2089.5
$ws.Range("M93").Value = -217.5
$ws.Range("H98").Value = 2688.158
$ws.Range("I98").Value = 4431.75
$ws.Range("J98").Value = 2223.2
$ws.Range("K98").Value = 13295.25
$ws.Range("L98").Value = 6669.599999999999
$ws.Range("M98").Value = -11797.25
$ws.Range("N98").Value = -9665.599999999999
$ws.Range("H105").Value = 38333.5
$ws.Range("I105").Value = 30001
$ws.Range("K105").Value = 90003
$ws.Range("M105").Value = -87382
$ws.Range("H110").Value = 1927
$ws.Range("I110").Value = 1927
$ws.Range("K110").Value = 5781
$ws.Range("M110").Value = -1691

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 30031.143
$ws.Range("J52").Value = 30031.143
$ws.Range("L52").Value = 30031.143
$ws.Range("N52").Value = -30549.143
$ws.Range("H113").Value = 724868.5
$ws.Range("J113").Value = 14194.8
$ws.Range("L113").Value = 14194.8
$ws.Range("N113").Value = -18534.8
$ws.Range("H134").Value = 63999.75
$ws.Range("J134").Value = 63999.75
$ws.Range("L134").Value = 191999.25
$ws.Range("N134").Value = -197069.25
$ws.Range("H136").Value = 55399
$ws.Range("J136").Value = 55399
$ws.Range("L136").Value = 166197
$ws.Range("N136").Value = -171297

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 10277
$ws.Range("I61").Value = 12036
$ws.Range("K61").Value = 12036
$ws.Range("M61").Value = -11834
$ws.Range("H113").Value = 10277
$ws.Range("I113").Value = 12036
$ws.Range("K113").Value = 12036
$ws.Range("M113").Value = -9866
$ws.Range("H133").Value = 56333.11
$ws.Range("J133").Value = 56333.11
$ws.Range("L133").Value = 56333.11
$ws.Range("N133").Value = -61393.11
$ws.Range("H136").Value = 136685.53
$ws.Range("I136").Value = 231072.36
$ws.Range("K136").Value = 693217.08
$ws.Range("M136").Value = -690667.08

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 641.3077
$ws.Range("I113").Value = 586.5909
$ws.Range("J113").Value = 942.25
$ws.Range("K113").Value = 1759.7727
$ws.Range("L113").Value = 2826.75
$ws.Range("M113").Value = 410.2273
$ws.Range("N113").Value = -7166.75
$ws.Range("H132").Value = 31854.195
$ws.Range("I132").Value = 4242.05
$ws.Range("K132").Value = 12726.15
$ws.Range("M132").Value = -10196.15
$ws.Range("H135").Value = 100049990
$ws.Range("J135").Value = 100049990
$ws.Range("N135").Value = -100060130
$ws.Range("L135").Value = 100049990
